$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Data Fields")

$ws.Range("D118").Value = 19091.93279569893
$ws.Range("D119").Value = 18825.34523809524
$ws.Range("D120").Value = 17232.53897849462
$ws.Range("D121").Value = 16683.79444444444
$ws.Range("D122").Value = 16196.16532258064
$ws.Range("D123").Value = 16970.33611111111
$ws.Range("D124").Value = 18394.55376344086
$ws.Range("D125").Value = 17483.94758064516
$ws.Range("D126").Value = 16126.75833333333
$ws.Range("D127").Value = 16221.31182795699
$ws.Range("D128").Value = 17593.7875
$ws.Range("D129").Value = 18630.69758064516
$ws.Range("D130").Value = 19906.20564516129
$ws.Range("D131").Value = 17990.8556547619
$ws.Range("D132").Value = 17432.80107526882
$ws.Range("D133").Value = 16188.03888888889
$ws.Range("D134").Value = 16319.80241935484
$ws.Range("D135").Value = 17168.31527777778
$ws.Range("D136").Value = 17132.45295698925
$ws.Range("D137").Value = 17084.39919354839
$ws.Range("D138").Value = 17149.18472222222
$ws.Range("D139").Value = 16431.7997311828
$ws.Range("D140").Value = 18045.04166666667
$ws.Range("D141").Value = 18366.11827956989
$ws.Range("D142").Value = 19906.20967741936
$ws.Range("B166").Value = 5946.377688172043
$ws.Range("C166").Value = 16271.64516129032
$ws.Range("B167").Value = 5755.407738095239
$ws.Range("C167").Value = 15788.09077380952
$ws.Range("T167").Value = 5946.377688172043
$ws.Range("B168").Value = 5770.697580645161
$ws.Range("C168").Value = 15580.03225806452
$ws.Range("T168").Value = 5755.407738095239
$ws.Range("B169").Value = 5120.201388888889
$ws.Range("C169").Value = 13597.32222222222
$ws.Range("T169").Value = 5770.697580645161
$ws.Range("B170").Value = 5182.504032258064
$ws.Range("C170").Value = 13646.29166666667
$ws.Range("T170").Value = 5120.201388888889
$ws.Range("B171").Value = 5679.984722222222
$ws.Range("C171").Value = 14799.18888888889
$ws.Range("T171").Value = 5182.504032258064
$ws.Range("B172").Value = 5951.712365591397
$ws.Range("C172").Value = 15572.87365591398
$ws.Range("T172").Value = 5679.984722222222
$ws.Range("B173").Value = 5852.674731182796
$ws.Range("C173").Value = 15259.55510752688
$ws.Range("T173").Value = 5951.712365591397
$ws.Range("B174").Value = 5727.665277777778
$ws.Range("C174").Value = 14833.28888888889
$ws.Range("T174").Value = 5852.674731182796
$ws.Range("B175").Value = 5207.611559139785
$ws.Range("C175").Value = 13894.97983870968
$ws.Range("T175").Value = 5727.665277777778
$ws.Range("B176").Value = 5524.329166666666
$ws.Range("C176").Value = 15214.7375
$ws.Range("T176").Value = 5207.611559139785
$ws.Range("B177").Value = 5692.653225806452
$ws.Range("C177").Value = 16493.22311827957
$ws.Range("T177").Value = 5524.329166666666
$ws.Range("B178").Value = 5916.83870967742
$ws.Range("C178").Value = 17072.60887096774
$ws.Range("T178").Value = 5692.653225806452
$ws.Range("T179").Value = 5916.83870967742
